$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 114, shifting existing rows 114:237 down to 115:238
$ws.Rows.Item(114).Insert()

# Populate the newly inserted row 114 with data (matches old row 114 values,
# except Fecha (D) and Volumen (J) which get new values)
$ws.Range("A114").Value = 3
$ws.Range("B114").Value = "Femacal de La Calera"
$ws.Range("C114").Value = "Coquimbo"
$ws.Range("D114").Value = 44539
$ws.Range("E114").Value = 5
$ws.Range("F114").Value = 100112012
$ws.Range("G114").Value = "Espinaca"
$ws.Range("H114").Value = "Sin especificar"
$ws.Range("I114").Value = "Primera"
$ws.Range("J114").Value = 160
$ws.Range("K114").Value = 2500
$ws.Range("L114").Value = 2500
$ws.Range("M114").Value = 2500
$ws.Range("N114").Value = '$/docena de atados (3 kilos)'
$ws.Range("O114").Value = "Provincia de Quillota"
$ws.Range("P114").Value = 833
$ws.Range("Q114").Value = 3
$ws.Range("R114").Value = "Hortaliza"
